# Automatische test-sync: 2025-06-18 11:00:10
#
# Appends two new "Klacht over levering" rows to the Logs sheet (rows 6 and 7),
# extends the conditional formatting ranges to cover the new rows, and
# refreshes the Dashboard summary (counts per Categorie, now including the
# extra "Klacht" entries, re-sorted to Klacht, Overig, Afmelding, Bestelling).

$wb = $excel.ActiveWorkbook

# --- Logs sheet: add the two new complaint rows -----------------------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A6").Value = "Klacht over levering"
$ws.Range("B6").Value = "mailmind.test@zohomail.eu"
$ws.Range("C6").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$ws.Range("D6").Value = "Klacht"
$ws.Range("F6").Value = "2025-06-18 10:30:11"
$ws.Range("G6").Value = "Nee"

$ws.Range("A7").Value = "Klacht over levering"
$ws.Range("B7").Value = "mailmind.test@zohomail.eu"
$ws.Range("C7").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$ws.Range("D7").Value = "Klacht"
$ws.Range("F7").Value = "2025-06-18 10:30:12"
$ws.Range("G7").Value = "Nee"

# --- Extend the conditional formatting ranges to include rows 6 and 7 ------
# The "Categorie" column formatting block (was D2:D5) -> D2:D7
$categorieRules = $ws.Range("D2:D5").FormatConditions
$categorieRules.Item(1).ModifyAppliesToRange($ws.Range("D2:D7"))

# The "Beantwoord" column formatting block (was G2:G5) -> G2:G7
$beantwoordRules = $ws.Range("G2:G5").FormatConditions
$beantwoordRules.Item(1).ModifyAppliesToRange($ws.Range("G2:G7"))

# --- Dashboard sheet: refresh the category counts ---------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A2").Value = "Klacht"
$dashboard.Range("B2").Value = 3
$dashboard.Range("A3").Value = "Overig"
$dashboard.Range("B3").Value = 1
$dashboard.Range("A4").Value = "Afmelding"
$dashboard.Range("B4").Value = 1
$dashboard.Range("A5").Value = "Bestelling"
$dashboard.Range("B5").Value = 1
